# Portfolio update for 2025-09-13: append a new row of price data,
# carrying forward the previous day's (2025-09-12) closing values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 29

# Column A holds the date as plain text (matches the existing rows, which
# store dates as text rather than Excel date serials). Forcing the cell to
# Text format before assignment prevents Excel from auto-converting the
# "yyyy-mm-dd" looking string into a date value/format, then ClearFormats
# removes the now-unneeded explicit formatting so the cell matches the
# unstyled look of the other data rows.
$dateCell = $ws.Range("A" + $newRow)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-09-13"
$dateCell.ClearFormats()

$ws.Range("B" + $newRow).Value = 57.11000061035156
$ws.Range("C" + $newRow).Value = 715.25
$ws.Range("D" + $newRow).Value = 321.3999938964844
